# Update the "想去人数" (interested-attendee count) figures on the
# 展览 (Exhibition) sheet and the 全部类型 (All Types) sheet, which
# duplicates the same events.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) — rows 2,3,5,6,7
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3366
$wsExhibit.Range("F3").Value = 17
$wsExhibit.Range("F5").Value = 1511
$wsExhibit.Range("F6").Value = 45
$wsExhibit.Range("F7").Value = 323

# Sheet "全部类型" (All Types) — rows 2,3,5,6,8 (row 7 here is a
# "演出" event that isn't present on the 展览 sheet, so the last
# updated row shifts from F7 to F8)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3366
$wsAll.Range("F3").Value = 17
$wsAll.Range("F5").Value = 1511
$wsAll.Range("F6").Value = 45
$wsAll.Range("F8").Value = 323
